$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.733.38"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.595.28"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "209.42"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "22.39"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "1.822.21"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "1.588.85"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "3.86"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "0.534"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").Value = "27.735.87"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "63.54"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "219.44"
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").Value = "0.0₃0698"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("D23").Value = "9.79"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").Value = "154.59"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").Value = "7.21"
$ws.Range("E26").Value = "  +4.33%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "15.16"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").Value = "0.105"
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "1.378.88"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("D34").Value = "2.98"
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("D36").Value = "0.980"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("D40").Value = "0.829"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "0.970"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("D43").Value = "64.62"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "1.74"
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "5.22"
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("D47").Value = "1.733.64"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "87.00"
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "0.0966"
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("D51").Value = "0.0496"
$ws.Range("E51").Value = "  -1.14%  "
